$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.266.29"
$ws.Range("E2").Value = "  +0.05%  "
$ws.Range("D3").Value = "1.899.36"
$ws.Range("E3").Value = "  -0.43%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").Value = "326.23"
$ws.Range("E5").Value = "  -0.54%  "
$ws.Range("E6").Value = "  -0.08%  "
$ws.Range("D7").Value = "0.4635"
$ws.Range("E7").Value = "  -0.60%  "
$ws.Range("E8").Value = "  -1.02%  "
$ws.Range("D9").Value = "0.07875"
$ws.Range("E9").Value = "  -1.26%  "
$ws.Range("D10").Value = "0.9887"
$ws.Range("E10").Value = "  -1.67%  "
$ws.Range("D11").Value = "21.83"
$ws.Range("E11").Value = "  -2.39%  "
$ws.Range("D12").Value = "1.904.42"
$ws.Range("E12").Value = "  -0.27%  "
$ws.Range("D13").Value = "7.062"
$ws.Range("E13").Value = "  -1.27%  "
$ws.Range("D14").Value = "5.733"
$ws.Range("D15").Value = "0.06976"
$ws.Range("E15").Value = "  +0.37%  "
$ws.Range("D16").Value = "88.30"
$ws.Range("E16").Value = "  -0.72%  "
$ws.Range("E17").Value = "  -0.10%  "
$ws.Range("D18").Value = "0.000009970"
$ws.Range("E18").Value = "  -1.62%  "
$ws.Range("D19").Value = "17.04"
$ws.Range("E19").Value = "  -1.00%  "
$ws.Range("E20").Value = "  -0.14%  "
$ws.Range("D21").Value = "29.272.78"
$ws.Range("E21").Value = "  +0.01%  "
$ws.Range("D22").Value = "5.298"
$ws.Range("E22").Value = "  -1.27%  "
$ws.Range("D23").Value = "11.09"
$ws.Range("E23").Value = "  +0.00%  "
$ws.Range("D25").Value = "156.03"
$ws.Range("E25").Value = "  -0.45%  "
$ws.Range("D26").Value = "19.38"
$ws.Range("E26").Value = "  -1.01%  "
$ws.Range("D27").Value = "6.033"
$ws.Range("E27").Value = "  +2.39%  "
$ws.Range("D28").Value = "118.51"
$ws.Range("E28").Value = "  -0.89%  "
$ws.Range("D29").Value = "1.890"
$ws.Range("E29").Value = "  -5.79%  "
$ws.Range("D30").Value = "0.09357"
$ws.Range("E30").Value = "  -0.92%  "
$ws.Range("D31").Value = "0.9032"
$ws.Range("E31").Value = "  -2.32%  "
$ws.Range("D32").Value = "5.251"
$ws.Range("E32").Value = "  -2.04%  "
$ws.Range("D33").Value = "1.323"
$ws.Range("E33").Value = "  -1.87%  "
$ws.Range("D34").Value = "3.214"
$ws.Range("E34").Value = "  -1.46%  "
$ws.Range("E35").Value = "  +1.10%  "
$ws.Range("D36").Value = "0.05775"
$ws.Range("E36").Value = "  -1.51%  "
$ws.Range("D37").Value = "0.02089"
$ws.Range("E37").Value = "  -1.06%  "
$ws.Range("D39").Value = "7.723"
$ws.Range("E39").Value = "  -3.45%  "
$ws.Range("D40").Value = "0.5704"
$ws.Range("E40").Value = "  -1.06%  "
$ws.Range("D41").Value = "0.1787"
$ws.Range("E41").Value = "  -1.50%  "
$ws.Range("D42").Value = "9.712"
$ws.Range("E42").Value = "  -3.29%  "
$ws.Range("D43").Value = "11.94"
$ws.Range("E43").Value = "  -0.74%  "
$ws.Range("D44").Value = "0.5356"
$ws.Range("E44").Value = "  -1.54%  "
$ws.Range("D45").Value = "2.182"
$ws.Range("E45").Value = "  -2.10%  "
$ws.Range("D46").Value = "0.07027"
$ws.Range("E46").Value = "  -0.88%  "
$ws.Range("D47").Value = "1.850"
$ws.Range("E47").Value = "  -2.05%  "
$ws.Range("D48").Value = "2.568"
$ws.Range("E48").Value = "  -0.47%  "
$ws.Range("D49").Value = "113.04"
$ws.Range("E49").Value = "  +0.84%  "
$ws.Range("D50").Value = "1.050"
$ws.Range("E50").Value = "  -2.03%  "
$ws.Range("B51").Value = "WOONetwork"
$ws.Range("C51").Value = "https://coinranking.com/coin/k-J3YwacF+woonetwork-woo"
$ws.Range("D51").Value = "0.2908"
$ws.Range("E51").Value = "  +0.00%  "
